$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: India (numbers refreshed) ---
$ws.Cells.Item(6,2).Value = 698817
$ws.Cells.Item(6,3).Value = 981
$ws.Cells.Item(6,4).Value = 424963
$ws.Cells.Item(6,5).Value = 254147
$ws.Cells.Item(6,7).Value = 7
$ws.Cells.Item(6,8).Value = 19707

# --- Rows 37/38: Ucrania overtakes Kazajistan in the ranking ---
# Row 37 becomes Ucrania with its refreshed totals
$ws.Cells.Item(37,1).Value = "Ucrania"
$ws.Cells.Item(37,2).Value = 49043
$ws.Cells.Item(37,3).Value = 543
$ws.Cells.Item(37,4).Value = 21703
$ws.Cells.Item(37,5).Value = 26078
$ws.Cells.Item(37,7).Value = 13
$ws.Cells.Item(37,8).Value = 1262

# Row 38 becomes Kazajistan, keeping its previous totals
$ws.Cells.Item(38,1).Value = "Kazajistan"
$ws.Cells.Item(38,2).Value = 48574
$ws.Cells.Item(38,3).Value = 1403
$ws.Cells.Item(38,4).Value = 27334
$ws.Cells.Item(38,5).Value = 21052
$ws.Cells.Item(38,7).Value = 0
$ws.Cells.Item(38,8).Value = 188

# --- Rows 52/53: Armenia overtakes Nigeria in the ranking ---
# Row 52 becomes Armenia with its refreshed totals
$ws.Cells.Item(52,1).Value = "Armenia"
$ws.Cells.Item(52,2).Value = 28936
$ws.Cells.Item(52,3).Value = 330
$ws.Cells.Item(52,4).Value = 16302
$ws.Cells.Item(52,5).Value = 12143
$ws.Cells.Item(52,7).Value = 7
$ws.Cells.Item(52,8).Value = 491

# Row 53 becomes Nigeria, keeping its previous totals
$ws.Cells.Item(53,1).Value = "Nigeria"
$ws.Cells.Item(53,2).Value = 28711
$ws.Cells.Item(53,3).Value = 0
$ws.Cells.Item(53,4).Value = 11665
$ws.Cells.Item(53,5).Value = 16401
$ws.Cells.Item(53,7).Value = 0
$ws.Cells.Item(53,8).Value = 645

# --- Row 77: El Salvador (numbers refreshed) ---
$ws.Cells.Item(77,4).Value = 4588
$ws.Cells.Item(77,5).Value = 2966
$ws.Cells.Item(77,7).Value = 6
$ws.Cells.Item(77,8).Value = 223

# --- Row 97: Hungria (numbers refreshed) ---
$ws.Cells.Item(97,2).Value = 4189
$ws.Cells.Item(97,3).Value = 6
$ws.Cells.Item(97,4).Value = 2860
$ws.Cells.Item(97,5).Value = 740

# --- Row 132: Letonia (numbers refreshed) ---
$ws.Cells.Item(132,2).Value = 1127
$ws.Cells.Item(132,3).Value = 3
$ws.Cells.Item(132,5).Value = 97

# --- Row 192: Islas Turcas y Caicos (numbers refreshed) ---
$ws.Cells.Item(192,2).Value = 48
$ws.Cells.Item(192,3).Value = 1
$ws.Cells.Item(192,5).Value = 35

# --- Timestamp footer update ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 6 de Julio de 2020 a las 09:25"
